$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the content of A1 to the new text for V2 of the worksheet.
$ws.Range("A1").Value = "This is V2 of the worksheet"

# Auto-fit column A so it's wide enough to display the new text.
$ws.Columns.Item(1).AutoFit() | Out-Null

# Move/park the selection on A2, matching the saved view state in the diff.
$ws.Range("A2").Select()
